$d = $word.ActiveDocument

# Locate the paragraph that contains the "Vorname" workshop-intro sentence
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Vorname*also took part in workshops*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $old = "{{ Vorname }} also took part in workshops with the following companies:"
    $new = "{%- if workshops|length > 0 %}{{ Vorname }} also took part in workshops with the following companies: {% endif %}"

    # Replace the text across the (currently split) runs, collapsing them into one run
    $target.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

    # Remove the paragraph's justified alignment (drop <w:jc w:val="both"/>)
    $target.Alignment = 0
}
